$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (row 26). All subsequent rows shift up by one.
$ws.Rows("26").Delete()

# After the shift, the old "SC 92" row is now row 27. Remove it too.
$ws.Rows("27").Delete()

# "SC 5" (now row 26) gains an "A" value that was previously missing.
$ws.Range("B26").Value = -20.2

# "SC 101" (now row 27) loses its "A" value, becoming blank again.
$ws.Range("B27").Value = ""
